$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.445.88"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +9.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.606.40"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +8.93%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.14"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +10.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9933"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.12%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.66%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3392"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +11.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.34"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +7.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.135"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +7.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07053"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +6.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +9.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.916"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +7.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.630"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +7.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.607.31"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +8.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001083"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +5.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9938"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +4.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06694"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +12.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "77.92"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +12.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.011"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +9.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.03"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +10.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.82"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +6.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.500.54"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +9.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.405"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +6.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.581"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +21.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "149.56"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +4.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.53"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +13.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.788.13"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +9.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "122.74"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +8.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.159"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +22.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.028"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9494"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +17.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.695"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +12.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08242"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.04"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +16.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.242"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +10.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.273"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.615"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +15.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06103"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02211"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +7.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2029"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +8.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9929"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5907"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +11.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.843"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +8.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.13"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +6.87%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5688"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +9.62%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.39"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +7.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.963"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +8.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06817"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +5.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.56"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +9.09%  "
